# Weekly fruit/vegetable price update: insert a new weekly record as row 12
# (pushing the existing rows 12-24 down to rows 13-25) for
# "Agrícola del Norte S.A. de Arica" - Perejil.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 12; Excel shifts rows 12:24 down to 13:25
# and copies row formatting from the row above (matches the style, incl. the
# date NumberFormat on column D, carried by the existing sheet).
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly observation.
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 44540
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112044
$ws.Range("G12").Value = "Perejil"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 950
$ws.Range("N12").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 475
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date/time number format.
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
